$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename columns to snake_case machine-readable names ---
$ws.Cells.Item(1, 1).Value = "mx_state"
$ws.Cells.Item(1, 2).Value = "mx_municipality"
$ws.Cells.Item(1, 3).Value = "n_matriculas"
$ws.Cells.Item(1, 4).Value = "pct_matriculas"

# --- Title-case connector words ("de", "del", "la", "las", "el", "los", "y") ---
# --- in state/municipality names (plus one stray-capitalization fix)      ---
$ws.Cells.Item(5, 2).Value = "Pabellón De Arteaga"
$ws.Cells.Item(22, 2).Value = "Amatenango De La Frontera"
$ws.Cells.Item(26, 2).Value = "Benemérito De Las Américas"
$ws.Cells.Item(32, 2).Value = "Comitán De Domínguez"
$ws.Cells.Item(49, 2).Value = "Mazapa De Madero"
$ws.Cells.Item(55, 2).Value = "Ocozocoautla De Espinosa"
$ws.Cells.Item(62, 2).Value = "Salto De Agua"
$ws.Cells.Item(63, 2).Value = "San Cristóbal De Las Casas"
$ws.Cells.Item(100, 1).Value = "Ciudad De México"
$ws.Cells.Item(104, 2).Value = "Cuajimalpa De Morelos"
$ws.Cells.Item(123, 2).Value = "Pánuco De Coronado"
$ws.Cells.Item(124, 2).Value = "San Juan Del Río"
$ws.Cells.Item(129, 1).Value = "Estado De México"
$ws.Cells.Item(129, 2).Value = "Acambay De Ruíz Castañeda"
$ws.Cells.Item(131, 2).Value = "Almoloya De Alquisiras"
$ws.Cells.Item(132, 2).Value = "Almoloya De Juárez"
$ws.Cells.Item(133, 2).Value = "Almoloya Del Río"
$ws.Cells.Item(138, 2).Value = "Atizapán De Zaragoza"
$ws.Cells.Item(145, 2).Value = "Chapa De Mota"
$ws.Cells.Item(148, 2).Value = "Coacalco De Berriozábal"
$ws.Cells.Item(155, 2).Value = "Ecatepec De Morelos"
$ws.Cells.Item(160, 2).Value = "Ixtapan De La Sal"
$ws.Cells.Item(161, 2).Value = "Ixtapan Del Oro"
$ws.Cells.Item(171, 2).Value = "Naucalpan De Juárez"
$ws.Cells.Item(180, 2).Value = "San Felipe Del Progreso"
$ws.Cells.Item(190, 2).Value = "Tenango Del Aire"
$ws.Cells.Item(191, 2).Value = "Tenango Del Valle"
$ws.Cells.Item(199, 2).Value = "Tlalnepantla De Baz"
$ws.Cells.Item(204, 2).Value = "Valle De Bravo"
$ws.Cells.Item(205, 2).Value = "Valle De Chalco Solidaridad"
$ws.Cells.Item(206, 2).Value = "Villa De Allende"
$ws.Cells.Item(207, 2).Value = "Villa Del Carbón"
$ws.Cells.Item(215, 2).Value = "Apaseo El Alto"
$ws.Cells.Item(216, 2).Value = "Apaseo El Grande"
$ws.Cells.Item(219, 2).Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Cells.Item(222, 2).Value = "Jaral Del Progreso"
$ws.Cells.Item(231, 2).Value = "San Francisco Del Rincón"
$ws.Cells.Item(232, 2).Value = "Silao De La Victoria"
$ws.Cells.Item(235, 2).Value = "Acapulco De Juárez"
$ws.Cells.Item(238, 2).Value = "Ajuchitlán Del Progreso"
$ws.Cells.Item(239, 2).Value = "Alcozauca De Guerrero"
$ws.Cells.Item(243, 2).Value = "Atenango Del Río"
$ws.Cells.Item(244, 2).Value = "Atlamajalcingo Del Monte"
$ws.Cells.Item(246, 2).Value = "Atoyac De Álvarez"
$ws.Cells.Item(247, 2).Value = "Ayutla De Los Libres"
$ws.Cells.Item(250, 2).Value = "Buenavista De Cuéllar"
$ws.Cells.Item(251, 2).Value = "Chilapa De Álvarez"
$ws.Cells.Item(252, 2).Value = "Chilpancingo De Los Bravo"
$ws.Cells.Item(253, 2).Value = "Coahuayutla De José María Izazaga"
$ws.Cells.Item(258, 2).Value = "Coyuca De Benítez"
$ws.Cells.Item(259, 2).Value = "Coyuca De Catalán"
$ws.Cells.Item(263, 2).Value = "Cuetzala Del Progreso"
$ws.Cells.Item(264, 2).Value = "Cutzamala De Pinzón"
$ws.Cells.Item(270, 2).Value = "Huitzuco De Los Figueroa"
$ws.Cells.Item(271, 2).Value = "Iguala De La Independencia"
$ws.Cells.Item(273, 2).Value = "Ixcateopan De Cuauhtémoc"
$ws.Cells.Item(274, 2).Value = "Zihuatanejo De Azueta"
$ws.Cells.Item(276, 2).Value = "La Unión De Isidoro Montes De Oca"
$ws.Cells.Item(279, 2).Value = "Mártir De Cuilapan"
$ws.Cells.Item(292, 2).Value = "Taxco De Alarcón"
$ws.Cells.Item(294, 2).Value = "Técpan De Galeana"
$ws.Cells.Item(296, 2).Value = "Tepecoacuilco De Trujano"
$ws.Cells.Item(298, 2).Value = "Tixtla De Guerrero"
$ws.Cells.Item(301, 2).Value = "Tlalixtaquilla De Maldonado"
$ws.Cells.Item(302, 2).Value = "Tlapa De Comonfort"
$ws.Cells.Item(314, 2).Value = "Agua Blanca De Iturbide"
$ws.Cells.Item(320, 2).Value = "Atotonilco De Tula"
$ws.Cells.Item(321, 2).Value = "Atotonilco El Grande"
$ws.Cells.Item(327, 2).Value = "Cuautepec De Hinojosa"
$ws.Cells.Item(331, 2).Value = "Huasca De Ocampo"
$ws.Cells.Item(334, 2).Value = "Huejutla De Reyes"
$ws.Cells.Item(337, 2).Value = "Jacala De Ledezma"
$ws.Cells.Item(343, 2).Value = "Mineral Del Chico"
$ws.Cells.Item(344, 2).Value = "Mineral Del Monte"
$ws.Cells.Item(345, 2).Value = "Mixquiahuala De Juárez"
$ws.Cells.Item(346, 2).Value = "Molango De Escamilla"
$ws.Cells.Item(348, 2).Value = "Nopala De Villagrán"
$ws.Cells.Item(349, 2).Value = "Pachuca De Soto"
$ws.Cells.Item(351, 2).Value = "Progreso De Obregón"
$ws.Cells.Item(356, 2).Value = "Santiago De Anaya"
$ws.Cells.Item(360, 2).Value = "Tenango De Doria"
$ws.Cells.Item(362, 2).Value = "Tepeji Del Río De Ocampo"
$ws.Cells.Item(364, 2).Value = "Tezontepec De Aldama"
$ws.Cells.Item(369, 2).Value = "Tula De Allende"
$ws.Cells.Item(370, 2).Value = "Tulancingo De Bravo"
$ws.Cells.Item(371, 2).Value = "Villa De Tezontepec"
$ws.Cells.Item(375, 2).Value = "Zacualtipán De Ángeles"
$ws.Cells.Item(378, 2).Value = "Ahualulco De Mercado"
$ws.Cells.Item(380, 2).Value = "Autlán De Navarro"
$ws.Cells.Item(388, 2).Value = "Lagos De Moreno"
$ws.Cells.Item(393, 2).Value = "San Cristóbal De La Barranca"
$ws.Cells.Item(394, 2).Value = "San Miguel El Alto"
$ws.Cells.Item(395, 2).Value = "Tamazula De Gordiano"
$ws.Cells.Item(397, 2).Value = "Tepatitlán De Morelos"
$ws.Cells.Item(399, 2).Value = "Tizapán El Alto"
$ws.Cells.Item(400, 2).Value = "Tlajomulco De Zúñiga"
$ws.Cells.Item(404, 2).Value = "Unión De San Antonio"
$ws.Cells.Item(405, 2).Value = "Valle De Juárez"
$ws.Cells.Item(406, 2).Value = "Yahualica De González Gallo"
$ws.Cells.Item(408, 2).Value = "Zapotlán El Grande"
$ws.Cells.Item(419, 2).Value = "Coalcomán De Vázquez Pallares"
$ws.Cells.Item(476, 2).Value = "Coatlán Del Río"
$ws.Cells.Item(487, 2).Value = "Puente De Ixtla"
$ws.Cells.Item(493, 2).Value = "Tetela Del Volcán"
$ws.Cells.Item(495, 2).Value = "Tlaltizapán De Zapata"
$ws.Cells.Item(503, 2).Value = "Zacualpan De Amilpas"
$ws.Cells.Item(505, 2).Value = "Amatlán De Cañas"
$ws.Cells.Item(518, 2).Value = "Montemorelos"
$ws.Cells.Item(523, 2).Value = "Acatlán De Pérez Figueroa"
$ws.Cells.Item(527, 2).Value = "Ayoquezco De Aldama"
$ws.Cells.Item(530, 2).Value = "Capulálpam De Méndez"
$ws.Cells.Item(532, 2).Value = "Chalcatongo De Hidalgo"
$ws.Cells.Item(533, 2).Value = "Ciénega De Zimatlán"
$ws.Cells.Item(535, 2).Value = "Coicoyán De Las Flores"
$ws.Cells.Item(536, 2).Value = "Constancia Del Rosario"
$ws.Cells.Item(539, 2).Value = "Cuyamecalco Villa De Zaragoza"
$ws.Cells.Item(540, 2).Value = "Fresnillo De Trujano"
$ws.Cells.Item(541, 2).Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Cells.Item(542, 2).Value = "Heroica Ciudad De Huajuapan De León"
$ws.Cells.Item(543, 2).Value = "Heroica Ciudad De Tlaxiaco"
$ws.Cells.Item(544, 2).Value = "Huautla De Jiménez"
$ws.Cells.Item(545, 2).Value = "Ixtlán De Juárez"
$ws.Cells.Item(546, 2).Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Cells.Item(553, 2).Value = "Magdalena Yodocono De Porfirio Díaz"
$ws.Cells.Item(555, 2).Value = "Mariscala De Juárez"
$ws.Cells.Item(557, 2).Value = "Mazatlán Villa De Flores"
$ws.Cells.Item(559, 2).Value = "Miahuatlán De Porfirio Díaz"
$ws.Cells.Item(560, 2).Value = "Nejapa De Madero"
$ws.Cells.Item(561, 2).Value = "Oaxaca De Juárez"
$ws.Cells.Item(562, 2).Value = "Ocotlán De Morelos"
$ws.Cells.Item(563, 2).Value = "Pinotepa De Don Luis"
$ws.Cells.Item(564, 2).Value = "Putla Villa De Guerrero"
$ws.Cells.Item(566, 2).Value = "Rojas De Cuauhtémoc"
$ws.Cells.Item(584, 2).Value = "San Baltazar Yatzachi El Bajo"
$ws.Cells.Item(617, 2).Value = "San Juan Del Estado"
$ws.Cells.Item(653, 2).Value = "San Miguel Del Puerto"
$ws.Cells.Item(654, 2).Value = "San Miguel El Grande"
$ws.Cells.Item(672, 2).Value = "San Pedro El Alto"
$ws.Cells.Item(684, 2).Value = "San Pedro Y San Pablo Teposcolula"
$ws.Cells.Item(685, 2).Value = "San Pedro Y San Pablo Tequixtepec"
$ws.Cells.Item(702, 2).Value = "Santa Cruz Tacache De Mina"
$ws.Cells.Item(705, 2).Value = "Santa Inés De Zaragoza"
$ws.Cells.Item(706, 2).Value = "Santa Inés Del Monte"
$ws.Cells.Item(707, 2).Value = "Santa Lucía Del Camino"
$ws.Cells.Item(721, 2).Value = "Santa María Jalapa Del Marqués"
$ws.Cells.Item(764, 2).Value = "Santo Domingo De Morelos"
$ws.Cells.Item(775, 2).Value = "Tamazulápam Del Espíritu Santo"
$ws.Cells.Item(776, 2).Value = "Tataltepec De Valdés"
$ws.Cells.Item(777, 2).Value = "Teotitlán De Flores Magón"
$ws.Cells.Item(778, 2).Value = "Tepelmeme Villa De Morelos"
$ws.Cells.Item(779, 2).Value = "Tezoatlán De Segura Y Luna"
$ws.Cells.Item(780, 2).Value = "Tlacolula De Matamoros"
$ws.Cells.Item(782, 2).Value = "Tlalixtac De Cabrera"
$ws.Cells.Item(783, 2).Value = "Totontepec Villa De Morelos"
$ws.Cells.Item(785, 2).Value = "Villa De Chilapa De Díaz"
$ws.Cells.Item(786, 2).Value = "Villa De Etla"
$ws.Cells.Item(787, 2).Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Cells.Item(788, 2).Value = "Villa De Zaachila"
$ws.Cells.Item(789, 2).Value = "Villa Sola De Vega"
$ws.Cells.Item(791, 2).Value = "Zimatlán De Álvarez"
$ws.Cells.Item(821, 2).Value = "Chalchicomula De Sesma"
$ws.Cells.Item(831, 2).Value = "Chila De La Sal"
$ws.Cells.Item(845, 2).Value = "Cuayuca De Andrade"
$ws.Cells.Item(846, 2).Value = "Cuetzalan Del Progreso"
$ws.Cells.Item(862, 2).Value = "Huehuetlán El Chico"
$ws.Cells.Item(863, 2).Value = "Huehuetlán El Grande"
$ws.Cells.Item(867, 2).Value = "Huitzilan De Serdán"
$ws.Cells.Item(869, 2).Value = "Ixcamilpa De Guerrero"
$ws.Cells.Item(873, 2).Value = "Izúcar De Matamoros"
$ws.Cells.Item(884, 2).Value = "Los Reyes De Juárez"
$ws.Cells.Item(885, 2).Value = "Mazapiltepec De Juárez"
$ws.Cells.Item(898, 2).Value = "Palmar De Bravo"
$ws.Cells.Item(908, 2).Value = "San Diego La Mesa Tochimiltzingo"
$ws.Cells.Item(923, 2).Value = "San Nicolás De Los Ranchos"
$ws.Cells.Item(927, 2).Value = "San Salvador El Seco"
$ws.Cells.Item(928, 2).Value = "San Salvador El Verde"
$ws.Cells.Item(935, 2).Value = "Tecali De Herrera"
$ws.Cells.Item(943, 2).Value = "Tepanco De López"
$ws.Cells.Item(944, 2).Value = "Tepango De Rodríguez"
$ws.Cells.Item(945, 2).Value = "Tepatlaxco De Hidalgo"
$ws.Cells.Item(951, 2).Value = "Tepexi De Rodríguez"
$ws.Cells.Item(953, 2).Value = "Tepeyahualco De Cuauhtémoc"
$ws.Cells.Item(954, 2).Value = "Tetela De Ocampo"
$ws.Cells.Item(959, 2).Value = "Tlacotepec De Benito Juárez"
$ws.Cells.Item(970, 2).Value = "Totoltepec De Guerrero"
$ws.Cells.Item(972, 2).Value = "Tuzamapan De Galeana"
$ws.Cells.Item(976, 2).Value = "Xayacatlán De Bravo"
$ws.Cells.Item(982, 2).Value = "Xochitlán De Vicente Suárez"
$ws.Cells.Item(989, 2).Value = "Zapotitlán De Méndez"
$ws.Cells.Item(997, 2).Value = "Amealco De Bonfil"
$ws.Cells.Item(999, 2).Value = "Cadereyta De Montes"
$ws.Cells.Item(1002, 2).Value = "Pinal De Amoles"
$ws.Cells.Item(1005, 2).Value = "San Juan Del Río"
$ws.Cells.Item(1011, 2).Value = "Armadillo De Los Infante"
$ws.Cells.Item(1012, 2).Value = "Ciudad Del Maíz"
$ws.Cells.Item(1021, 2).Value = "Santa María Del Río"
$ws.Cells.Item(1027, 2).Value = "Villa De Guadalupe"
$ws.Cells.Item(1028, 2).Value = "Villa De Reyes"
$ws.Cells.Item(1061, 2).Value = "Jalpa De Méndez"
$ws.Cells.Item(1082, 2).Value = "Acuamanala De Miguel Hidalgo"
$ws.Cells.Item(1084, 2).Value = "Apetatitlán De Antonio Carvajal"
$ws.Cells.Item(1089, 2).Value = "Contla De Juan Cuamatzi"
$ws.Cells.Item(1096, 2).Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Cells.Item(1099, 2).Value = "Mazatecochco De José María Morelos"
$ws.Cells.Item(1100, 2).Value = "Muñoz De Domingo Arenas"
$ws.Cells.Item(1101, 2).Value = "Nanacamilpa De Mariano Arista"
$ws.Cells.Item(1104, 2).Value = "Papalotla De Xicohténcatl"
$ws.Cells.Item(1110, 2).Value = "San Pablo Del Monte"
$ws.Cells.Item(1111, 2).Value = "Sanctórum De Lázaro Cárdenas"
$ws.Cells.Item(1120, 2).Value = "Tepetitla De Lardizábal"
$ws.Cells.Item(1123, 2).Value = "Tetla De La Solidaridad"
$ws.Cells.Item(1140, 2).Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Cells.Item(1143, 2).Value = "Amatlán De Los Reyes"
$ws.Cells.Item(1150, 2).Value = "Camarón De Tejeda"
$ws.Cells.Item(1153, 2).Value = "Castillo De Teayo"
$ws.Cells.Item(1160, 2).Value = "Chinampa De Gorostiza"
$ws.Cells.Item(1169, 2).Value = "Cosamaloapan De Carpio"
$ws.Cells.Item(1181, 2).Value = "Hueyapan De Ocampo"
$ws.Cells.Item(1182, 2).Value = "Huiloapan De Cuauhtémoc"
$ws.Cells.Item(1183, 2).Value = "Ignacio De La Llave"
$ws.Cells.Item(1187, 2).Value = "Ixhuatlán De Madero"
$ws.Cells.Item(1188, 2).Value = "Ixhuatlán Del Café"
$ws.Cells.Item(1189, 2).Value = "Ixhuatlán Del Sureste"
$ws.Cells.Item(1197, 2).Value = "Juchique De Ferrer"
$ws.Cells.Item(1202, 2).Value = "Las Vigas De Ramírez"
$ws.Cells.Item(1203, 2).Value = "Lerdo De Tejada"
$ws.Cells.Item(1206, 2).Value = "Martínez De La Torre"
$ws.Cells.Item(1210, 2).Value = "Mixtla De Altamirano"
$ws.Cells.Item(1216, 2).Value = "Paso De Ovejas"
$ws.Cells.Item(1217, 2).Value = "Paso Del Macho"
$ws.Cells.Item(1220, 2).Value = "Poza Rica De Hidalgo"
$ws.Cells.Item(1227, 2).Value = "Sayula De Alemán"
$ws.Cells.Item(1230, 2).Value = "Soledad De Doblado"
$ws.Cells.Item(1233, 2).Value = "Tatahuicapan De Juárez"
$ws.Cells.Item(1259, 2).Value = "Vega De Alatorre"
$ws.Cells.Item(1267, 2).Value = "Zontecomatlán De López Y Fuentes"
$ws.Cells.Item(1268, 2).Value = "Zozocolco De Hidalgo"

# --- Remove trailing footnote/metadata rows (1286:1290) and shrink the  ---
# --- used range / dimension from D1290 down to D1284                     ---
$ws.Range("A1286:A1290").EntireRow.Delete() | Out-Null

Write-Output $ws.UsedRange.Address()
